$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NSE:ASHOKAMET"
$ws.Range("C2").Value = "NSE:AHLADA"
$ws.Range("E2").Value = "NSE:ALKEM"
$ws.Range("F2").ClearContents()
$ws.Range("B3").Value = "NSE:BALPHARMA"
$ws.Range("C3").Value = "NSE:ALKALI"
$ws.Range("E3").Value = "NSE:BAJAJFINSV"
$ws.Range("B4").Value = "NSE:ICEMAKE"
$ws.Range("C4").Value = "NSE:AMBER"
$ws.Range("E4").Value = "NSE:BAJFINANCE"
$ws.Range("C5").Value = "NSE:ANANTRAJ"
$ws.Range("E5").Value = "NSE:BALKRISIND"
$ws.Range("C6").Value = "NSE:APOLLO"
$ws.Range("E6").Value = "NSE:INFY"
$ws.Range("C7").Value = "NSE:ATFL"
$ws.Range("E7").Value = "NSE:LTIM"
$ws.Range("C8").Value = "NSE:BALKRISHNA"
$ws.Range("E8").Value = "NSE:MANAPPURAM"
$ws.Range("C9").Value = "NSE:BALMLAWRIE"
$ws.Range("E9").Value = "NSE:NTPC"
$ws.Range("C10").Value = "NSE:BORORENEW"
$ws.Range("E10").Value = "NSE:RELIANCE"
$ws.Range("C11").Value = "NSE:CAREERP"
$ws.Range("C12").Value = "NSE:CCL"
$ws.Range("C13").Value = "NSE:COASTCORP"
$ws.Range("C14").Value = "NSE:CYIENT"
$ws.Range("C15").Value = "NSE:DCW"
$ws.Range("C16").Value = "NSE:DECCANCE"
$ws.Range("C17").Value = "NSE:DISHTV"
$ws.Range("C18").Value = "NSE:EKC"
$ws.Range("C19").Value = "NSE:EMUDHRA"
$ws.Range("C20").Value = "NSE:EVERESTIND"
$ws.Range("C21").Value = "NSE:FACT"
$ws.Range("C22").Value = "NSE:FOSECOIND"
$ws.Range("C23").Value = "NSE:GFLLIMITED"
$ws.Range("C24").Value = "NSE:GINNIFILA"
$ws.Range("C25").Value = "NSE:GOACARBON"
$ws.Range("C26").Value = "NSE:GODREJPROP"
$ws.Range("C27").Value = "NSE:GOLDIAM"
$ws.Range("C28").Value = "NSE:GRAVITA"
$ws.Range("C29").Value = "NSE:IMAGICAA"
$ws.Range("C30").Value = "NSE:JBMA"
$ws.Range("C31").Value = "NSE:JINDRILL"
$ws.Range("C32").Value = "NSE:JMFINANCIL"
$ws.Range("C33").Value = "NSE:JSL"
$ws.Range("C34").Value = "NSE:KECL"
$ws.Range("C36").Value = "NSE:LICI"
$ws.Range("C37").Value = "NSE:MAHASTEEL"
$ws.Range("C38").Value = "NSE:MICEL"
$ws.Range("C39").Value = "NSE:MOLDTECH"
$ws.Range("C40").Value = "NSE:NAVNETEDUL"
$ws.Range("C41").Value = "NSE:NELCO"
$ws.Range("C42").Value = "NSE:NEWGEN"
$ws.Range("C43").Value = "NSE:NLCINDIA"
$ws.Range("C44").Value = "NSE:NSIL"
$ws.Range("C45").Value = "NSE:NYKAA"
$ws.Range("C46").Value = "NSE:OMINFRAL"
$ws.Range("C47").Value = "NSE:PATINTLOG"
$ws.Range("C48").Value = "NSE:PFS"
$ws.Range("C49").Value = "NSE:PGHL"
$ws.Range("C50").Value = "NSE:PILITA"
$ws.Range("C52").Value = "NSE:RAILTEL"
$ws.Range("C53").Value = "NSE:SAKSOFT"

$ws.Rows("54:56").Delete()
